# Add team record (Wins/Losses/Ties) columns to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they match the rest of the header row (bold, bordered,
# centered) instead of getting default formatting.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player (rows 2-54) gets the same team record.
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 57
    $ws.Cells.Item($r, 31).Value = 105
    $ws.Cells.Item($r, 32).Value = 0
}
